$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp in the header/footer cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 21:07"

# 2) Refresh the daily country statistics. The sheet is kept sorted by
#    column B ("Casos totales") descending, so a handful of countries
#    change rank (and therefore row) as part of this data refresh -
#    for those rows the country name (column A) is updated too.
$countryData = @(
    @(4, "Estados Unidos", 7119713, 21776, 4373627, 2540019, 0, 597, 206067),
    @(5, "India", 5727750, 87254, 4667060, 969527, 0, 1142, 91163),
    @(11, "España", 693556, 11289, 0, 0, 0, 130, 31034),
    @(14, "Francia", 481141, 13072, 93538, 356144, 0, 43, 31459),
    @(25, "Alemania", 278631, 1455, 247900, 21226, 0, 14, 9505),
    @(29, "Canada", 147515, 852, 127403, 10870, 0, 8, 9242),
    @(101, "Guayana Francesa", 9762, 24, 9431, 266, 0, 0, 65),
    @(116, "Cabo Verde", 5412, 75, 4837, 521, 0, 2, 54),
    @(117, "Republica de Yibuti", 5407, 0, 5338, 8, 0, 0, 61),
    @(118, "Jamaica", 5395, 125, 1444, 3875, 0, 1, 76),
    @(130, "Angola", 4363, 127, 1473, 2731, 0, 4, 159),
    @(135, "Aruba", 3721, 56, 2501, 1195, 0, 0, 25),
    @(136, "Gambia", 3542, 2, 2011, 1421, 0, 0, 110),
    @(137, "Mayotte", 3541, 0, 2964, 537, 0, 0, 40),
    @(139, "Reunion", 3501, 86, 2482, 1004, 0, 0, 15),
    @(140, "Bahamas", 3467, 0, 1871, 1519, 0, 0, 77),
    @(141, "Somalia", 3465, 0, 2877, 490, 0, 0, 98),
    @(146, "Sudan del Sur", 2664, 4, 1290, 1325, 0, 0, 49),
    @(159, "Belice", 1669, 34, 1004, 644, 0, 0, 21),
    @(162, "Lesoto", 1507, 83, 766, 706, 0, 2, 35),
    @(186, "Curazao", 301, 10, 104, 196, 0, 0, 1),
    @(190, "Monaco", 199, 2, 165, 33, 0, 0, 1),
    @(214, "Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
    @(215, "Montserrat", 13, 0, 12, 0, 0, 0, 1)
)

foreach ($r in $countryData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($rowNum, 2 + $c).Value = $r[2 + $c]
    }
}